$d = $word.ActiveDocument

# 1. "Adresse des Creditors" -> "Adresse D C"
$d.Content.Find.Execute('„Adresse des Creditors“', $false, $false, $false, $false, $false, $true, 1, $false, "„Adresse D C“", 2) | Out-Null

# 2. "heutiges Datum" (straight quotes) -> "heutiges D" (curly quotes)
$d.Content.Find.Execute('"heutiges Datum"', $false, $false, $false, $false, $false, $true, 1, $false, "„heutiges D“", 2) | Out-Null

# 3. "Aktenzeichen des Creditors" + trailing space -> "Aktenzeichen D C" + trailing space
$d.Content.Find.Execute('„Aktenzeichen des Creditors“ ', $false, $false, $false, $false, $false, $true, 1, $false, "„Aktenzeichen D C“ ", 2) | Out-Null

# 4. merge "„Name“" + "," runs
$d.Content.Find.Execute('„Name“,', $false, $false, $false, $false, $false, $true, 1, $false, "„Name“,", 2) | Out-Null

# 5. merge "„Geburtstag“" + ", wohnhaft" runs
$d.Content.Find.Execute('„Geburtstag“, wohnhaft', $false, $false, $false, $false, $false, $true, 1, $false, "„Geburtstag“, wohnhaft", 2) | Out-Null

# 6. merge "„Adresse“" + long sentence, dropping the trailing "Ordnungsgemäße..." sentence
$d.Content.Find.Execute('„Adresse“, wird von uns bei der Durchführung eines außergerichtlichen Einigungsversuchs im Rahmen des Verbraucherinsolvenzverfahrens gemäß § 305, Abs. 1 Nr. 1 InsO vertreten. Ordnungsgemäße Bevollmächtigung wird anwaltlich versichert.', $false, $false, $false, $false, $false, $true, 1, $false, "„Adresse“, wird von uns bei der Durchführung eines außergerichtlichen Einigungsversuchs im Rahmen des Verbraucherinsolvenzverfahrens gemäß § 305, Abs. 1 Nr. 1 InsO vertreten. ", 2) | Out-Null

# 7. merge "„Name“" + " strebt eine Schuldenbereinigung..." runs
$d.Content.Find.Execute('„Name“ strebt eine Schuldenbereinigung auf der Grundlage der ', $false, $false, $false, $false, $false, $true, 1, $false, "„Name“ strebt eine Schuldenbereinigung auf der Grundlage der ", 2) | Out-Null

# 8. merge "Da-ten" hyphenation + proofErr removal
$d.Content.Find.Execute('der Angelegenheit zu veranlassen. Erfahrungsgemäß dauert es einige Zeit, bis uns alle relevanten Da-ten vorliegen. Sobald dies der Fall ist, kommen wir unaufgefordert wieder auf Sie zu.', $false, $false, $false, $false, $false, $true, 1, $false, "der Angelegenheit zu veranlassen. Erfahrungsgemäß dauert es einige Zeit, bis uns alle relevanten Daten vorliegen. Sobald dies der Fall ist, kommen wir unaufgefordert wieder auf Sie zu.", 2) | Out-Null

# 9. fix "gebe-ten." -> "gebeten."
$d.Content.Find.Execute('gebe-ten.', $false, $false, $false, $false, $false, $true, 1, $false, "gebeten.", 2) | Out-Null

# 10. append " Thomas Scuric" after the final "Rechtsanwalt" signature line (last paragraph of the body)
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)
$lastPara.Range.InsertAfter(" Thomas Scuric")

Write-Output "done"
